# Work Time Table.xlsx edit
# - Fill in the next work day (A4) using the same date-format style as the
#   rows above it (A2/A3), and set its value to 31-May-2020 (serial 43982).
# - Move the active selection to B10 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from A2 onto A4, then set A4's value.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = 43982

# Update the saved selection / active cell.
$ws.Range("B10").Select() | Out-Null
